$wb = $excel.ActiveWorkbook

# ALC!row18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 933.63635
$ws.Range("I18").Value = 692.381
$ws.Range("K18").Value = 692.381
$ws.Range("M18").Value = -408.381

# ALC!row40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2560.652
$ws.Range("I40").Value = 2615.3845
$ws.Range("J40").Value = 2489.5
$ws.Range("K40").Value = 2615.3845
$ws.Range("L40").Value = 2489.5
$ws.Range("M40").Value = -2440.3845
$ws.Range("N40").Value = -2839.5

# ALC!row43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1869.2858
$ws.Range("I43").Value = 1841.375
$ws.Range("J43").Value = 1906.5
$ws.Range("K43").Value = 1841.375
$ws.Range("L43").Value = 1906.5
$ws.Range("M43").Value = -1772.375
$ws.Range("N43").Value = -2044.5

# ALC!row64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 70399.13
$ws.Range("J64").Value = 3999.0715
$ws.Range("L64").Value = 3999.0715
$ws.Range("N64").Value = -4495.0715

# ALC!row67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 70399.13
$ws.Range("J67").Value = 3999.0715
$ws.Range("L67").Value = 3999.0715
$ws.Range("N67").Value = -5715.0715

# ALC!row70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2171
$ws.Range("J70").Value = 1459.4
$ws.Range("L70").Value = 4378.200000000001
$ws.Range("N70").Value = -4918.200000000001

# ALC!row73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 2171
$ws.Range("J73").Value = 1459.4
$ws.Range("L73").Value = 4378.200000000001
$ws.Range("N73").Value = -6250.200000000001

# ALC!row74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# ALC!row77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# ALC!row129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 305037.47
$ws.Range("J129").Value = 357150.56
$ws.Range("L129").Value = 1071451.68
$ws.Range("N129").Value = -1081451.68

# ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22967.375
$ws.Range("I32").Value = 3987.2463
$ws.Range("K32").Value = 3987.2463
$ws.Range("M32").Value = -3700.2463

# ARM!row63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3217.5
$ws.Range("I63").Value = 1605
$ws.Range("J63").Value = 3540
$ws.Range("K63").Value = 1605
$ws.Range("L63").Value = 3540
$ws.Range("M63").Value = -919
$ws.Range("N63").Value = -4912

# ARM!row66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3217.5
$ws.Range("I66").Value = 1605
$ws.Range("J66").Value = 3540
$ws.Range("K66").Value = 8025
$ws.Range("L66").Value = 17700
$ws.Range("M66").Value = -4593
$ws.Range("N66").Value = -24564

# BSM!row35
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 17766.666
$ws.Range("J35").Value = 17766.666
$ws.Range("L35").Value = 17766.666
$ws.Range("N35").Value = -18386.666

# BSM!row36
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 7387.4443
$ws.Range("I36").Value = 8305.75
$ws.Range("J36").Value = 41
$ws.Range("K36").Value = 8305.75
$ws.Range("L36").Value = 41
$ws.Range("M36").Value = -7771.75
$ws.Range("N36").Value = -1109

# BSM!row45
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H45").Value = 12688.333
$ws.Range("J45").Value = 12688.333
$ws.Range("L45").Value = 12688.333
$ws.Range("N45").Value = -14304.333

# BSM!row101
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H101").Value = 25577
$ws.Range("J101").Value = 25577
$ws.Range("L101").Value = 25577
$ws.Range("N101").Value = -32067

# CRP!row41
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 8463.75
$ws.Range("I41").Value = 3683.3333
$ws.Range("J41").Value = 11332
$ws.Range("K41").Value = 3683.3333
$ws.Range("L41").Value = 11332
$ws.Range("M41").Value = -3255.3333
$ws.Range("N41").Value = -12188

# CRP!row51
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 7329.077
$ws.Range("I51").Value = 500
$ws.Range("J51").Value = 7898.1665
$ws.Range("K51").Value = 500
$ws.Range("L51").Value = 7898.1665
$ws.Range("M51").Value = 236
$ws.Range("N51").Value = -9370.166499999999

# CRP!row59
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 22265.715
$ws.Range("J59").Value = 24643.334
$ws.Range("L59").Value = 24643.334
$ws.Range("N59").Value = -26933.334

# CRP!row60
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 9467.883
$ws.Range("I60").Value = 7600
$ws.Range("J60").Value = 10246.167
$ws.Range("K60").Value = 7600
$ws.Range("L60").Value = 10246.167
$ws.Range("M60").Value = -7089
$ws.Range("N60").Value = -11268.167

# CRP!row61
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 7329.077
$ws.Range("I61").Value = 500
$ws.Range("J61").Value = 7898.1665
$ws.Range("K61").Value = 500
$ws.Range("L61").Value = 7898.1665
$ws.Range("M61").Value = -152
$ws.Range("N61").Value = -8594.166499999999

# CRP!row122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 616.25
$ws.Range("I122").Value = 581
$ws.Range("J122").Value = 651.5
$ws.Range("K122").Value = 1743
$ws.Range("L122").Value = 1954.5
$ws.Range("M122").Value = 707
$ws.Range("N122").Value = -6854.5

# CRP!row132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3471.524
$ws.Range("I132").Value = 3323.2942
$ws.Range("K132").Value = 9969.882599999999
$ws.Range("M132").Value = -7439.882599999999

# CUL!row11
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 261.1
$ws.Range("I11").Value = 141.83333
$ws.Range("J11").Value = 440
$ws.Range("K11").Value = 425.49999
$ws.Range("L11").Value = 1320
$ws.Range("M11").Value = -285.49999
$ws.Range("N11").Value = -1600

# LTW!row46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1125085.5
$ws.Range("I46").Value = 369.5
$ws.Range("J46").Value = 1446432.9
$ws.Range("K46").Value = 369.5
$ws.Range("L46").Value = 1446432.9
$ws.Range("M46").Value = -181.5
$ws.Range("N46").Value = -1446808.9

# LTW!row55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 355471.4
$ws.Range("I55").Value = 631508.9399999999
$ws.Range("J55").Value = 566
$ws.Range("K55").Value = 631508.9399999999
$ws.Range("L55").Value = 566
$ws.Range("M55").Value = -631335.9399999999
$ws.Range("N55").Value = -912

# LTW!row68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2848.9443
$ws.Range("I68").Value = 1456
$ws.Range("J68").Value = 3735.3635
$ws.Range("K68").Value = 1456
$ws.Range("L68").Value = 3735.3635
$ws.Range("M68").Value = -707
$ws.Range("N68").Value = -5233.363499999999

# LTW!row71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2848.9443
$ws.Range("I71").Value = 1456
$ws.Range("J71").Value = 3735.3635
$ws.Range("K71").Value = 7280
$ws.Range("L71").Value = 18676.8175
$ws.Range("M71").Value = -3536
$ws.Range("N71").Value = -26164.8175

# LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2342.6567
$ws.Range("I132").Value = 2364.1372
$ws.Range("J132").Value = 2274.1875
$ws.Range("K132").Value = 7092.4116
$ws.Range("L132").Value = 6822.5625
$ws.Range("M132").Value = -4562.4116
$ws.Range("N132").Value = -11882.5625

# WVR!row96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 100001590
$ws.Range("I96").Value = 200001860
$ws.Range("J96").Value = 1320
$ws.Range("K96").Value = 200001860
$ws.Range("L96").Value = 1320
$ws.Range("M96").Value = -200000487
$ws.Range("N96").Value = -4066
